$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 286.89655
$ws.Range("I80").Value = 284.66666
$ws.Range("J80").Value = 288.47058
$ws.Range("K80").Value = 853.9999799999999
$ws.Range("L80").Value = 865.41174
$ws.Range("M80").Value = 144.0000200000001
$ws.Range("N80").Value = -2861.41174
$ws.Range("H83").Value = 286.89655
$ws.Range("I83").Value = 284.66666
$ws.Range("J83").Value = 288.47058
$ws.Range("K83").Value = 2561.99994
$ws.Range("L83").Value = 2596.23522
$ws.Range("M83").Value = 2430.00006
$ws.Range("N83").Value = -12580.23522
$ws.Range("H132").Value = 2061.6667
$ws.Range("I132").Value = 2182
$ws.Range("J132").Value = 1460
$ws.Range("K132").Value = 6546
$ws.Range("L132").Value = 4380
$ws.Range("M132").Value = -4016
$ws.Range("N132").Value = -9440
$ws.Range("H137").Value = 10419320
$ws.Range("I137").Value = 27780144
$ws.Range("J137").Value = 2824.6667
$ws.Range("K137").Value = 83340432
$ws.Range("L137").Value = 8474.000100000001
$ws.Range("M137").Value = -83337882
$ws.Range("N137").Value = -13574.0001
$ws.Range("H138").Value = 2959.627
$ws.Range("I138").Value = 1161.2727
$ws.Range("J138").Value = 3838.8223
$ws.Range("K138").Value = 3483.8181
$ws.Range("L138").Value = 11516.4669
$ws.Range("M138").Value = 1656.1819
$ws.Range("N138").Value = -21796.4669
$ws.Range("H139").Value = 10883.846
$ws.Range("J139").Value = 10883.846
$ws.Range("L139").Value = 10883.846
$ws.Range("N139").Value = -21163.846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1998.2354
$ws.Range("I2").Value = 2096.7693
$ws.Range("J2").Value = 1678
$ws.Range("K2").Value = 2096.7693
$ws.Range("L2").Value = 1678
$ws.Range("M2").Value = -1983.7693
$ws.Range("N2").Value = -1904
$ws.Range("H33").Value = 6026
$ws.Range("I33").Value = 6026
$ws.Range("K33").Value = 6026
$ws.Range("M33").Value = -5697
$ws.Range("H61").Value = 2334.8928
$ws.Range("I61").Value = 1703.1052
$ws.Range("J61").Value = 3668.6667
$ws.Range("K61").Value = 1703.1052
$ws.Range("L61").Value = 3668.6667
$ws.Range("M61").Value = -1491.1052
$ws.Range("N61").Value = -4092.6667
$ws.Range("H116").Value = 1998.2354
$ws.Range("I116").Value = 2096.7693
$ws.Range("J116").Value = 1678
$ws.Range("K116").Value = 2096.7693
$ws.Range("L116").Value = 1678
$ws.Range("M116").Value = 197.2307000000001
$ws.Range("N116").Value = -6266
$ws.Range("H132").Value = 2582.2856
$ws.Range("I132").Value = 2199.8518
$ws.Range("J132").Value = 3873
$ws.Range("K132").Value = 6599.555399999999
$ws.Range("L132").Value = 11619
$ws.Range("M132").Value = -4069.555399999999
$ws.Range("N132").Value = -16679
$ws.Range("H136").Value = 2334.8928
$ws.Range("I136").Value = 1703.1052
$ws.Range("J136").Value = 3668.6667
$ws.Range("K136").Value = 5109.3156
$ws.Range("L136").Value = 11006.0001
$ws.Range("M136").Value = -2559.3156
$ws.Range("N136").Value = -16106.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1998.2354
$ws.Range("I3").Value = 2096.7693
$ws.Range("J3").Value = 1678
$ws.Range("K3").Value = 2096.7693
$ws.Range("L3").Value = 1678
$ws.Range("M3").Value = -1982.7693
$ws.Range("N3").Value = -1906
$ws.Range("H134").Value = 3059.2
$ws.Range("I134").Value = 2810.8
$ws.Range("J134").Value = 3804.4
$ws.Range("K134").Value = 8432.400000000001
$ws.Range("L134").Value = 11413.2
$ws.Range("M134").Value = -5897.400000000001
$ws.Range("N134").Value = -16483.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1431140.8
$ws.Range("I62").Value = 1669247.5
$ws.Range("K62").Value = 1669247.5
$ws.Range("M62").Value = -1668623.5
$ws.Range("H65").Value = 1431140.8
$ws.Range("I65").Value = 1669247.5
$ws.Range("K65").Value = 8346237.5
$ws.Range("M65").Value = -8343117.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 19.5
$ws.Range("I8").Value = 19.5
$ws.Range("K8").Value = 58.5
$ws.Range("M8").Value = 80.5
$ws.Range("H68").Value = 1114.0526
$ws.Range("I68").Value = 1106.6842
$ws.Range("J68").Value = 1121.421
$ws.Range("K68").Value = 3320.0526
$ws.Range("L68").Value = 3364.263
$ws.Range("M68").Value = -2509.0526
$ws.Range("N68").Value = -4986.263
$ws.Range("H71").Value = 1114.0526
$ws.Range("I71").Value = 1106.6842
$ws.Range("J71").Value = 1121.421
$ws.Range("K71").Value = 9960.157799999999
$ws.Range("L71").Value = 10092.789
$ws.Range("M71").Value = -5904.157799999999
$ws.Range("N71").Value = -18204.789
$ws.Range("H131").Value = 1106.3334
$ws.Range("J131").Value = 1129.6428
$ws.Range("L131").Value = 3388.9284
$ws.Range("N131").Value = -13468.9284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1675.409
$ws.Range("I102").Value = 1697.8948
$ws.Range("K102").Value = 1697.8948
$ws.Range("M102").Value = -75.89480000000003
$ws.Range("H126").Value = 3699.8572
$ws.Range("I126").Value = 3144.2222
$ws.Range("K126").Value = 9432.6666
$ws.Range("M126").Value = -6962.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 701.9286
$ws.Range("I16").Value = 735.4167
$ws.Range("J16").Value = 501
$ws.Range("K16").Value = 735.4167
$ws.Range("L16").Value = 501
$ws.Range("M16").Value = -565.4167
$ws.Range("N16").Value = -841
$ws.Range("H32").Value = 1799.6666
$ws.Range("I32").Value = 1799.6666
$ws.Range("K32").Value = 1799.6666
$ws.Range("M32").Value = -1482.6666
$ws.Range("H68").Value = 180229.34
$ws.Range("I68").Value = 527972.9399999999
$ws.Range("J68").Value = 1658.2972
$ws.Range("K68").Value = 527972.9399999999
$ws.Range("L68").Value = 1658.2972
$ws.Range("M68").Value = -527223.9399999999
$ws.Range("N68").Value = -3156.2972
$ws.Range("H71").Value = 180229.34
$ws.Range("I71").Value = 527972.9399999999
$ws.Range("J71").Value = 1658.2972
$ws.Range("K71").Value = 2639864.7
$ws.Range("L71").Value = 8291.486000000001
$ws.Range("M71").Value = -2636120.7
$ws.Range("N71").Value = -15779.486

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 7000
$ws.Range("J42").Value = 7000
$ws.Range("L42").Value = 7000
$ws.Range("N42").Value = -7756
$ws.Range("H122").Value = 3004.5264
$ws.Range("I122").Value = 2403.25
$ws.Range("J122").Value = 3441.818
$ws.Range("K122").Value = 7209.75
$ws.Range("L122").Value = 10325.454
$ws.Range("M122").Value = -4759.75
$ws.Range("N122").Value = -15225.454
